# Refresh cryptocurrency price/volume data to the latest scrape (GitHub Actions bot update)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new text would NOT be misread as a number - assign directly.
$textUpdates = @{
    'D2' = '66.099.46'
    'E2' = '  +0.07%  '
    'D3' = '3.258.49'
    'E3' = '  +2.50%  '
    'E4' = '  +0.05%  '
    'E5' = '  +0.51%  '
    'E6' = '  +2.04%  '
    'E7' = '  +0.10%  '
    'D8' = '3.258.37'
    'E8' = '  +2.41%  '
    'E9' = '  -0.76%  '
    'E10' = '  +1.38%  '
    'E11' = '  +1.29%  '
    'E12' = '  -2.76%  '
    'E13' = '  +1.92%  '
    'E14' = '  +0.85%  '
    'D15' = '3.795.33'
    'E15' = '  +2.54%  '
    'D16' = '66.160.47'
    'E16' = '  +0.08%  '
    'D17' = '3.263.83'
    'E17' = '  +2.63%  '
    'E18' = '  -0.52%  '
    'E19' = '  +1.16%  '
    'E20' = '  -1.57%  '
    'E21' = '  +0.23%  '
    'E22' = '  +2.96%  '
    'E23' = '  +0.76%  '
    'E24' = '  -1.03%  '
    'E25' = '  +2.56%  '
    'E26' = '  -0.07%  '
    'E27' = '  +1.10%  '
    'E28' = '  +0.34%  '
    'E29' = '  +46.79%  '
    'E30' = '  -0.47%  '
    'E31' = '  +0.10%  '
    'E32' = '  -5.90%  '
    'E33' = '  -0.30%  '
    'E34' = '  -0.08%  '
    'E35' = '  -2.67%  '
    'E36' = '  -0.77%  '
    'E37' = '  +21.22%  '
    'B38' = 'OKB'
    'C38' = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
    'E38' = '  +0.53%  '
    'B39' = 'PEPE'
    'C39' = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
    'D39' = '0.0₃0789'
    'E39' = '  +11.58%  '
    'E40' = '  -2.51%  '
    'E41' = '  +1.49%  '
    'E42' = '  +1.58%  '
    'E43' = '  +0.30%  '
    'E44' = '  +4.75%  '
    'D45' = '2.993.75'
    'E45' = '  +5.81%  '
    'E46' = '  -2.37%  '
    'E47' = '  +3.07%  '
    'E48' = '  +4.48%  '
    'E49' = '  +2.18%  '
    'E51' = '  -1.83%  '
}

# Cells whose new text looks like a plain number (e.g. "1.00", "0.540").
# These must keep their original "Price" column look (trailing zeros, etc.)
# exactly as text, so force the cell to Text format before writing the value.
$numericLookingUpdates = @{
    'D4' = '1.00'
    'D5' = '606.99'
    'D6' = '157.50'
    'D9' = '0.540'
    'D10' = '0.161'
    'D11' = '5.71'
    'D12' = '0.495'
    'D13' = '0.0000270'
    'D14' = '38.58'
    'D18' = '7.31'
    'D20' = '499.19'
    'D21' = '15.31'
    'D22' = '0.749'
    'D23' = '8.05'
    'D24' = '14.63'
    'D25' = '86.45'
    'D26' = '0.999'
    'D27' = '3.02'
    'D28' = '9.16'
    'D29' = '0.132'
    'D30' = '2.36'
    'D31' = '7.03'
    'D32' = '2.84'
    'D33' = '27.80'
    'D35' = '1.14'
    'D36' = '6.42'
    'D37' = '3.48'
    'D38' = '55.55'
    'D40' = '498.35'
    'D41' = '0.0422'
    'D42' = '0.129'
    'D43' = '8.77'
    'D44' = '2.57'
    'D46' = '0.291'
    'D47' = '28.79'
    'D48' = '2.46'
    'D51' = '121.20'
}

foreach ($key in $textUpdates.Keys) {
    $ws.Range($key).Value = $textUpdates[$key]
}

foreach ($key in $numericLookingUpdates.Keys) {
    $cell = $ws.Range($key)
    $cell.NumberFormat = "@"
    $cell.Value = $numericLookingUpdates[$key]
}
